$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32, pushing existing rows 32-110 down to 33-111.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Range("A32").Value = 5
$ws.Range("B32").Value = "Macroferia Regional de Talca"
$ws.Range("C32").Value = "Maule"
$ws.Range("D32").Value = 45148
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = 100112040
$ws.Range("G32").Value = "Cilantro"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 150
$ws.Range("K32").Value = 8000
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 8000
$ws.Range("N32").Value = "$/caja 36 atados"
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 222
$ws.Range("Q32").Value = 36
$ws.Range("R32").Value = "Hortaliza"
